$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 1
$ws.Range("D4").Value = "Atrasos no desenvolvimento da Plataforma de Automação de Testes de Software"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "O atraso pode comprometer a eficiência das equipes de desenvolvimento e reduzir a competitividade no mercado."
$ws.Range("H4").Value = "Mitigar"
$ws.Range("I4").Value = "Implementar metodologias ágeis e monitoramento frequente dos prazos"
$ws.Rows.Item(4).RowHeight = 46.8

$ws.Range("B5").Value = 2
$ws.Range("D5").Value = "Baixa adesão à Plataforma de Treinamento para Equipes de Desenvolvimento"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "As equipes podem não utilizar adequadamente, resultando em baixo ROI."
$ws.Range("H5").Value = "Mitigar"
$ws.Range("I5").Value = "Programa de capacitação e incentivo para o uso da plataforma"
$ws.Rows.Item(5).RowHeight = 46.8

$ws.Range("B6").Value = 3
$ws.Range("D6").Value = "`tMudanças constantes nos requisitos da Plataforma de Gestão de Requisitos"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "Mudanças frequentes podem gerar retrabalho, custos e prazos aumentados."
$ws.Range("H6").Value = "Mitigar"
$ws.Range("I6").Value = "Definir escopo claro e processos robustos de gerenciamento de mudanças"
$ws.Rows.Item(6).RowHeight = 46.8

$ws.Range("B7").Value = 4
$ws.Range("D7").Value = "`tFalta de inovação no Projeto de Pesquisa e Melhoria Contínua dos Processos"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = "Falta de melhorias pode resultar na perda de competitividade e estagnação."
$ws.Range("H7").Value = "Mitigar"
$ws.Range("I7").Value = "Investir em pesquisa contínua e parcerias com universidades ou startups"
$ws.Rows.Item(7).RowHeight = 46.8

$ws.Range("B8").Value = 5
$ws.Range("D8").Value = "Falha de integração com sistemas legados no Laboratório de Inovação para Engenharia de Software"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = "A incompatibilidade tecnológica pode atrasar a implementação e causar custos adicionais."
$ws.Range("H8").Value = "Mitigar"
$ws.Range("I8").Value = "`tPlanejar fases de teste rigorosas e utilização de APIs bem documentadas"
$ws.Rows.Item(8).RowHeight = 46.8

$ws.Range("B9").Value = 6
$ws.Range("D9").Value = "Sobrecarga de recursos em múltiplos projetos simultâneos"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = "A equipe pode ficar sobrecarregada, resultando em atrasos e baixa qualidade."
$ws.Range("H9").Value = "Mitigar"
$ws.Range("I9").Value = "Alocar recursos com base na prioridade e monitorar a carga de trabalho"
$ws.Rows.Item(9).RowHeight = 46.8

$ws.Range("B10").Value = 7
$ws.Range("D10").Value = "Falta de feedback contínuo para o Projeto de Pesquisa e Melhoria Contínua dos Processos"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = "Sem feedback adequado, as melhorias podem ser mal direcionadas."
$ws.Range("H10").Value = "Mitigar"
$ws.Range("I10").Value = "Implementar ciclos curtos de feedback com stakeholders e usuários"
$ws.Rows.Item(10).RowHeight = 46.8

$ws.Range("B11").Value = 8
$ws.Range("D11").Value = "Custos elevados no desenvolvimento da Plataforma de Gestão de Requisitos"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = "O aumento nos custos pode comprometer o ROI e a viabilidade do projeto."
$ws.Range("H11").Value = "Mitigar"
$ws.Range("I11").Value = "Acompanhamento contínuo do orçamento e revisão das funcionalidades"
$ws.Rows.Item(11).RowHeight = 46.8

$ws.Range("H9").Select()
